$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2799.923
$ws.Range("I113").Value = 2720
$ws.Range("J113").Value = 2849.875
$ws.Range("K113").Value = 2720
$ws.Range("L113").Value = 2849.875
$ws.Range("M113").Value = 534
$ws.Range("N113").Value = -9357.875
$ws.Range("H116").Value = 3098
$ws.Range("I116").Value = 2974.4443
$ws.Range("J116").Value = 3468.6667
$ws.Range("K116").Value = 2974.4443
$ws.Range("L116").Value = 3468.6667
$ws.Range("M116").Value = 467.5556999999999
$ws.Range("N116").Value = -10352.6667
$ws.Range("H132").Value = 3357.0715
$ws.Range("I132").Value = 3038.423
$ws.Range("J132").Value = 7499.5
$ws.Range("K132").Value = 9115.269
$ws.Range("L132").Value = 22498.5
$ws.Range("M132").Value = -6585.269
$ws.Range("N132").Value = -27558.5
$ws.Range("H138").Value = 2114.9429
$ws.Range("I138").Value = 1815.5186
$ws.Range("K138").Value = 5446.5558
$ws.Range("M138").Value = -306.5558000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 564452
$ws.Range("I32").Value = 653871.1
$ws.Range("J32").Value = 18002
$ws.Range("K32").Value = 653871.1
$ws.Range("L32").Value = 18002
$ws.Range("M32").Value = -653584.1
$ws.Range("N32").Value = -18576
$ws.Range("H74").Value = 978
$ws.Range("I74").Value = 658
$ws.Range("J74").Value = 2115.7778
$ws.Range("K74").Value = 658
$ws.Range("L74").Value = 2115.7778
$ws.Range("M74").Value = 216
$ws.Range("N74").Value = -3863.7778
$ws.Range("H77").Value = 978
$ws.Range("I77").Value = 658
$ws.Range("J77").Value = 2115.7778
$ws.Range("K77").Value = 3290
$ws.Range("L77").Value = 10578.889
$ws.Range("M77").Value = 1078
$ws.Range("N77").Value = -19314.889
$ws.Range("H132").Value = 3306.7017
$ws.Range("I132").Value = 2075.8948
$ws.Range("J132").Value = 5768.316
$ws.Range("K132").Value = 6227.6844
$ws.Range("L132").Value = 17304.948
$ws.Range("M132").Value = -3697.6844
$ws.Range("N132").Value = -22364.948

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 25642798
$ws.Range("I20").Value = 1747
$ws.Range("J20").Value = 71430380
$ws.Range("K20").Value = 1747
$ws.Range("L20").Value = 71430380
$ws.Range("M20").Value = -1500
$ws.Range("N20").Value = -71430874
$ws.Range("H107").Value = 501500
$ws.Range("I107").Value = 1000000
$ws.Range("K107").Value = 1000000
$ws.Range("M107").Value = -998080
$ws.Range("H134").Value = 1557.0878
$ws.Range("I134").Value = 1205.591
$ws.Range("J134").Value = 2746.7693
$ws.Range("K134").Value = 3616.773
$ws.Range("L134").Value = 8240.3079
$ws.Range("M134").Value = -1081.773
$ws.Range("N134").Value = -13310.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 843.8333
$ws.Range("I16").Value = 687.5
$ws.Range("K16").Value = 687.5
$ws.Range("M16").Value = -400.5
$ws.Range("H113").Value = 843.8333
$ws.Range("I113").Value = 687.5
$ws.Range("K113").Value = 687.5
$ws.Range("M113").Value = 1482.5
$ws.Range("H122").Value = 1928.3043
$ws.Range("I122").Value = 1764
$ws.Range("J122").Value = 1952.95
$ws.Range("K122").Value = 5292
$ws.Range("L122").Value = 5858.85
$ws.Range("M122").Value = -2842
$ws.Range("N122").Value = -10758.85

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 856.6
$ws.Range("I5").Value = 587.4
$ws.Range("J5").Value = 1395
$ws.Range("K5").Value = 1762.2
$ws.Range("L5").Value = 4185
$ws.Range("M5").Value = -1650.2
$ws.Range("N5").Value = -4409
$ws.Range("H9").Value = 45486.7
$ws.Range("I9").Value = 333
$ws.Range("J9").Value = 47863.21
$ws.Range("K9").Value = 999
$ws.Range("L9").Value = 143589.63
$ws.Range("M9").Value = -775
$ws.Range("N9").Value = -144037.63
$ws.Range("H69").Value = 2046.0435
$ws.Range("I69").Value = 425.42856
$ws.Range("J69").Value = 2755.0625
$ws.Range("K69").Value = 1276.28568
$ws.Range("L69").Value = 8265.1875
$ws.Range("M69").Value = -465.28568
$ws.Range("N69").Value = -9887.1875
$ws.Range("H72").Value = 2046.0435
$ws.Range("I72").Value = 425.42856
$ws.Range("J72").Value = 2755.0625
$ws.Range("K72").Value = 3828.85704
$ws.Range("L72").Value = 24795.5625
$ws.Range("M72").Value = 227.1429600000001
$ws.Range("N72").Value = -32907.5625
$ws.Range("H107").Value = 62500412
$ws.Range("I107").Value = 467.66666
$ws.Range("J107").Value = 100000380
$ws.Range("K107").Value = 1402.99998
$ws.Range("L107").Value = 300001140
$ws.Range("M107").Value = 517.0000199999999
$ws.Range("N107").Value = -300004980
$ws.Range("H131").Value = 1030.6666
$ws.Range("J131").Value = 1138.1923
$ws.Range("L131").Value = 3414.5769
$ws.Range("N131").Value = -13494.5769
$ws.Range("H132").Value = 2414.6487
$ws.Range("J132").Value = 2567.56
$ws.Range("L132").Value = 23108.04
$ws.Range("N132").Value = -28168.04
$ws.Range("H135").Value = 856.6
$ws.Range("I135").Value = 587.4
$ws.Range("J135").Value = 1395
$ws.Range("K135").Value = 5286.599999999999
$ws.Range("L135").Value = 12555
$ws.Range("M135").Value = -2751.599999999999
$ws.Range("N135").Value = -17625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 26600
$ws.Range("J92").Value = 26600
$ws.Range("L92").Value = 26600
$ws.Range("N92").Value = -30344
$ws.Range("H132").Value = 6563.625
$ws.Range("I132").Value = 7499.6665
$ws.Range("J132").Value = 6002
$ws.Range("K132").Value = 22498.9995
$ws.Range("L132").Value = 18006
$ws.Range("M132").Value = -19968.9995
$ws.Range("N132").Value = -23066

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7770.95
$ws.Range("I22").Value = 2357.1428
$ws.Range("J22").Value = 10686.077
$ws.Range("K22").Value = 2357.1428
$ws.Range("L22").Value = 10686.077
$ws.Range("M22").Value = -2062.1428
$ws.Range("N22").Value = -11276.077
$ws.Range("H27").Value = 7770.95
$ws.Range("I27").Value = 2357.1428
$ws.Range("J27").Value = 10686.077
$ws.Range("K27").Value = 2357.1428
$ws.Range("L27").Value = 10686.077
$ws.Range("M27").Value = -2250.1428
$ws.Range("N27").Value = -10900.077
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("H132").Value = 2996.1765
$ws.Range("I132").Value = 2089.8
$ws.Range("J132").Value = 4291
$ws.Range("K132").Value = 6269.400000000001
$ws.Range("L132").Value = 12873
$ws.Range("M132").Value = -3739.400000000001
$ws.Range("N132").Value = -17933
$ws.Range("H133").Value = 34860.8
$ws.Range("J133").Value = 34860.8
$ws.Range("L133").Value = 34860.8
$ws.Range("N133").Value = -39920.8
$ws.Range("H136").Value = 1071.1666
$ws.Range("I136").Value = 934.89655
$ws.Range("J136").Value = 1635.7142
$ws.Range("K136").Value = 2804.68965
$ws.Range("L136").Value = 4907.142599999999
$ws.Range("M136").Value = -254.6896500000003
$ws.Range("N136").Value = -10007.1426

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 47500
$ws.Range("J92").Value = 47500
$ws.Range("L92").Value = 47500
$ws.Range("N92").Value = -52492
